# Auto-generated edit script applying targeted cell value updates
# per the provided diff (Jogos_da_Semana_FlashScore_2025-02-07.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G10").Value = 3.2
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 2.45
$ws.Range("L10").Value = 3.4
$ws.Range("Y10").Value = 7
$ws.Range("AA10").Value = 13
$ws.Range("AD10").Value = 51
$ws.Range("AK10").Value = 10
$ws.Range("AM10").Value = 23
$ws.Range("AR10").Value = 4.8
$ws.Range("H11").Value = 2.88
$ws.Range("I11").Value = 2.7
$ws.Range("K11").Value = 1.83
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("Q11").Value = 2.88
$ws.Range("R11").Value = 1.4
$ws.Range("AA11").Value = 12
$ws.Range("AE11").Value = 5.5
$ws.Range("AF11").Value = 6
$ws.Range("AL11").Value = 11
$ws.Range("AM11").Value = 26
$ws.Range("AO11").Value = 41
$ws.Range("AP11").Value = 2.2
$ws.Range("AQ11").Value = 1.68
$ws.Range("AR11").Value = 5.2
$ws.Range("AS11").Value = 1.17
$ws.Range("O14").Value = 1.91
$ws.Range("P14").Value = 1.8
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 1.06
$ws.Range("G22").Value = 1.5
$ws.Range("H22").Value = 3.9
$ws.Range("O22").Value = 1.4
$ws.Range("P22").Value = 2.75
$ws.Range("Q22").Value = 2.25
$ws.Range("R22").Value = 1.62
$ws.Range("S22").Value = 4.33
$ws.Range("T22").Value = 1.2
$ws.Range("W22").Value = 2.5
$ws.Range("X22").Value = 1.5
$ws.Range("Y22").Value = 5
$ws.Range("AB22").Value = 10
$ws.Range("AE22").Value = 7.5
$ws.Range("AG22").Value = 26
$ws.Range("AL22").Value = 21
$ws.Range("AR22").Value = 3.4
$ws.Range("G45").Value = 1.48
$ws.Range("I45").Value = 6
$ws.Range("J45").Value = 2
$ws.Range("O45").Value = 1.17
$ws.Range("P45").Value = 5
$ws.Range("Q45").Value = 1.57
$ws.Range("R45").Value = 2.35
$ws.Range("S45").Value = 2.38
$ws.Range("T45").Value = 1.53
$ws.Range("U45").Value = 1.29
$ws.Range("V45").Value = 3.5
$ws.Range("W45").Value = 1.7
$ws.Range("X45").Value = 2.05
$ws.Range("Y45").Value = 8.5
$ws.Range("Z45").Value = 8
$ws.Range("AB45").Value = 11
$ws.Range("AD45").Value = 21
$ws.Range("AH45").Value = 41
$ws.Range("AI45").Value = 201
$ws.Range("AJ45").Value = 19
$ws.Range("AL45").Value = 17
$ws.Range("AR45").Value = 1.98
$ws.Range("AS45").Value = 1.88
$ws.Range("G46").Value = 2.55
$ws.Range("I46").Value = 2.7
$ws.Range("L46").Value = 3.4
$ws.Range("O46").Value = 1.33
$ws.Range("P46").Value = 3.25
$ws.Range("R46").Value = 1.75
$ws.Range("Y46").Value = 8
$ws.Range("Z46").Value = 12
$ws.Range("AB46").Value = 23
$ws.Range("AM46").Value = 29
$ws.Range("AO46").Value = 34
$ws.Range("G47").Value = 1.8
$ws.Range("I47").Value = 3.75
$ws.Range("J47").Value = 2.35
$ws.Range("L47").Value = 4.1
$ws.Range("O47").Value = 1.23
$ws.Range("Q47").Value = 1.7
$ws.Range("R47").Value = 1.93
$ws.Range("S47").Value = 2.6
$ws.Range("T47").Value = 1.38
$ws.Range("X47").Value = 2
$ws.Range("Y47").Value = 8
$ws.Range("Z47").Value = 9.25
$ws.Range("AA47").Value = 8.25
$ws.Range("AB47").Value = 15
$ws.Range("AC47").Value = 13.5
$ws.Range("AG47").Value = 14
$ws.Range("AH47").Value = 55
$ws.Range("AJ47").Value = 12.5
$ws.Range("AK47").Value = 22
$ws.Range("AL47").Value = 12.5
$ws.Range("AM47").Value = 55
$ws.Range("AN47").Value = 32
$ws.Range("AO47").Value = 35
